$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-10 04:48:19"
$ws.Range("I2").Value = "10.3 mm"
$ws.Range("E3").Value = "2026-02-10 04:48:21"
$ws.Range("I3").Value = "6.1 mm"
$ws.Range("L3").Value = "42.8 km/h - 243º 4:12 TU"
$ws.Range("E4").Value = "2026-02-10 04:48:23"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "86%"
$ws.Range("H4").NumberFormat = "General"
$ws.Range("J4").Value = "1004.8 hPa"
$ws.Range("N4").Value = "7.5 °C 4:25 TU"
$ws.Range("O4").Value = "9.0 °C"
$ws.Range("E5").Value = "2026-02-10 04:48:25"
$ws.Range("I5").Value = "9.6 mm"
$ws.Range("L5").Value = "34.2 km/h - 325º 4:00 TU"
$ws.Range("O5").Value = "-0.3 °C"
$ws.Range("E6").Value = "2026-02-10 04:48:28"
$ws.Range("J6").Value = "1004.9 hPa"
$ws.Range("M6").Value = "8.0 °C 4:07 TU"
$ws.Range("N6").Value = "6.7 °C 4:29 TU"
$ws.Range("E7").Value = "2026-02-10 04:48:30"
$ws.Range("N7").Value = "11.6 °C 4:09 TU"
$ws.Range("E8").Value = "2026-02-10 04:48:32"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "95%"
$ws.Range("H8").NumberFormat = "General"
$ws.Range("J8").Value = "1005.1 hPa"
$ws.Range("E9").Value = "2026-02-10 04:48:35"
$ws.Range("N9").Value = "5.3 °C 4:22 TU"
$ws.Range("E10").Value = "2026-02-10 04:48:37"
$ws.Range("E11").Value = "2026-02-10 04:48:39"
$ws.Range("E12").Value = "2026-02-10 04:48:42"
$ws.Range("N12").Value = "5.8 °C 4:29 TU"
$ws.Range("O12").Value = "7.2 °C"
$ws.Range("E13").Value = "2026-02-10 04:48:44"
$ws.Range("I13").Value = "1.5 mm"
$ws.Range("N13").Value = "2.4 °C 4:00 TU"
$ws.Range("O13").Value = "2.6 °C"
$ws.Range("E14").Value = "2026-02-10 04:48:46"
$ws.Range("N14").Value = "8.2 °C 4:01 TU"
$ws.Range("O14").Value = "9.6 °C"
$ws.Range("E15").Value = "2026-02-10 04:48:49"
$ws.Range("N15").Value = "4.3 °C 4:23 TU"
$ws.Range("O15").Value = "6.9 °C"
$ws.Range("E16").Value = "2026-02-10 04:48:51"
$ws.Range("I16").Value = "8.6 mm"
$ws.Range("E17").Value = "2026-02-10 04:48:53"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "97%"
$ws.Range("H17").NumberFormat = "General"
$ws.Range("M17").Value = "4.5 °C 4:21 TU"
$ws.Range("O17").Value = "2.0 °C"
$ws.Range("E18").Value = "2026-02-10 04:48:56"
$ws.Range("N18").Value = "6.6 °C 4:16 TU"
$ws.Range("O18").Value = "7.8 °C"
$ws.Range("E19").Value = "2026-02-10 04:48:58"
$ws.Range("N19").Value = "3.6 °C 4:16 TU"
$ws.Range("E20").Value = "2026-02-10 04:49:00"
$ws.Range("I20").Value = "2.0 mm"
$ws.Range("O20").Value = "-1.3 °C"
$ws.Range("E21").Value = "2026-02-10 04:49:03"
$ws.Range("I21").Value = "2.3 mm"
$ws.Range("J21").Value = "1007.4 hPa"
$ws.Range("E22").Value = "2026-02-10 04:49:05"
$ws.Range("G22").Value = "123 cm"
$ws.Range("M22").Value = "-1.5 °C 4:22 TU"
$ws.Range("O22").Value = "-2.0 °C"
$ws.Range("E23").Value = "2026-02-10 04:49:07"
$ws.Range("I23").Value = "7.0 mm"
$ws.Range("O23").Value = "-0.5 °C"
$ws.Range("E24").Value = "2026-02-10 04:49:10"
$ws.Range("I24").Value = "0.9 mm"
$ws.Range("O24").Value = "8.3 °C"
$ws.Range("E25").Value = "2026-02-10 04:49:12"
$ws.Range("G25").Value = "116 cm"
$ws.Range("I25").Value = "4.4 mm"
$ws.Range("E26").Value = "2026-02-10 04:49:15"
$ws.Range("E27").Value = "2026-02-10 04:49:17"
$ws.Range("O27").Value = "-0.5 °C"
$ws.Range("E28").Value = "2026-02-10 04:49:20"
$ws.Range("N28").Value = "4.2 °C 4:16 TU"
$ws.Range("O28").Value = "5.4 °C"
$ws.Range("E29").Value = "2026-02-10 04:49:22"
$ws.Range("N29").Value = "7.5 °C 4:29 TU"
$ws.Range("O29").Value = "9.3 °C"
$ws.Range("E30").Value = "2026-02-10 04:49:24"
$ws.Range("L30").Value = "13.3 km/h - 23º 4:12 TU"
$ws.Range("N30").Value = "7.0 °C 4:03 TU"
$ws.Range("E31").Value = "2026-02-10 04:49:27"
$ws.Range("J31").Value = "1004.4 hPa"
$ws.Range("N31").Value = "7.8 °C 4:10 TU"
$ws.Range("O31").Value = "9.0 °C"
$ws.Range("E32").Value = "2026-02-10 04:49:29"
$ws.Range("O32").Value = "7.6 °C"
$ws.Range("E33").Value = "2026-02-10 04:49:32"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "97%"
$ws.Range("H33").NumberFormat = "General"
$ws.Range("I33").Value = "1.9 mm"
$ws.Range("E34").Value = "2026-02-10 04:49:34"
$ws.Range("M34").Value = "3.5 °C 4:29 TU"
$ws.Range("O34").Value = "2.5 °C"
$ws.Range("E35").Value = "2026-02-10 04:49:36"
$ws.Range("N35").Value = "10.3 °C 4:28 TU"
$ws.Range("O35").Value = "10.5 °C"
$ws.Range("E36").Value = "2026-02-10 04:49:39"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "92%"
$ws.Range("H36").NumberFormat = "General"
$ws.Range("J36").Value = "1004.9 hPa"
$ws.Range("N36").Value = "7.7 °C 4:29 TU"
$ws.Range("E37").Value = "2026-02-10 04:49:41"
$ws.Range("O37").Value = "4.0 °C"
$ws.Range("E38").Value = "2026-02-10 04:49:44"
$ws.Range("E39").Value = "2026-02-10 04:49:46"
$ws.Range("I39").Value = "1.5 mm"
$ws.Range("E40").Value = "2026-02-10 04:49:48"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "98%"
$ws.Range("H40").NumberFormat = "General"
$ws.Range("I40").Value = "2.9 mm"
$ws.Range("J40").Value = "1008.1 hPa"
$ws.Range("N40").Value = "4.3 °C 4:00 TU"
$ws.Range("O40").Value = "4.8 °C"
$ws.Range("E41").Value = "2026-02-10 04:49:51"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "98%"
$ws.Range("H41").NumberFormat = "General"
$ws.Range("N41").Value = "8.3 °C 4:29 TU"
$ws.Range("O41").Value = "10.0 °C"
$ws.Range("E42").Value = "2026-02-10 04:49:53"
$ws.Range("N42").Value = "7.2 °C 4:26 TU"
$ws.Range("O42").Value = "8.3 °C"
$ws.Range("E43").Value = "2026-02-10 04:49:55"
$ws.Range("M43").Value = "6.5 °C 4:02 TU"
$ws.Range("E44").Value = "2026-02-10 04:49:58"
$ws.Range("I44").Value = "6.2 mm"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-10 04:50:00"
$ws.Range("I45").Value = "11.9 mm"
$ws.Range("J45").Value = "1008.0 hPa"
$ws.Range("L45").Value = "15.5 km/h - 115º 4:15 TU"
$ws.Range("M45").Value = "4.0 °C 4:20 TU"
$ws.Range("O45").Value = "3.1 °C"
$ws.Range("E46").Value = "2026-02-10 04:50:02"
$ws.Range("N46").Value = "8.9 °C 4:03 TU"
$ws.Range("O46").Value = "9.8 °C"
